$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.610.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.495.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.61%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.494.32'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.577'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.73%  '

$ws.Range("E11").Value = '  +5.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.435'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.098.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.50%  '

$ws.Range("E14").Value = '  -0.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.629.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.45%  '

$ws.Range("E17").Value = '  +3.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.522.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '388.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.30%  '

$ws.Range("E24").Value = '  -0.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.531'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.49%  '

$ws.Range("E26").Value = '  +5.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.51%  '

$ws.Range("E28").Value = '  +3.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.15%  '

$ws.Range("E30").Value = '  +6.71%  '

$ws.Range("E31").Value = '  +6.80%  '

$ws.Range("E32").Value = '  +2.61%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.52'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.34%  '

$ws.Range("E34").Value = '  +6.09%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("E36").Value = '  +7.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.26%  '

$ws.Range("E38").Value = '  +6.04%  '

$ws.Range("E39").Value = '  +5.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.49%  '

$ws.Range("E41").Value = '  +2.75%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.49%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.30'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.814.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.11%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.64'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '357.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0310'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.88%  '

$ws.Range("E50").Value = '  +4.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.65%  '
